$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 119, shifting the existing rows 119-144 down to 120-145
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new weekly record
$ws.Cells.Item(119, 1).Value = 10
$ws.Cells.Item(119, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(119, 3).Value = "La Araucanía"
$ws.Cells.Item(119, 4).Value = Get-Date -Year 2023 -Month 9 -Day 4 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(119, 5).Value = 9
$ws.Cells.Item(119, 6).Value = 100112010
$ws.Cells.Item(119, 7).Value = "Achicoria"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 120
$ws.Cells.Item(119, 11).Value = 10000
$ws.Cells.Item(119, 12).Value = 10000
$ws.Cells.Item(119, 13).Value = 10000
$ws.Cells.Item(119, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(119, 15).Value = "Región Metropolitana"
$ws.Cells.Item(119, 16).Value = 556
$ws.Cells.Item(119, 17).Value = 18
$ws.Cells.Item(119, 18).Value = "Hortaliza"
